$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three highlighted "self-loop" cells in the automaton transition table
# (E4, F6, G7) are being cleared back to the default "-" / unmarked style,
# matching the style already used by plain cells such as D4, E6 and F7.

# E4: q2 (highlighted) -> "-" (plain), matching D4's style
$ws.Range("E4").Value = "-"
$ws.Range("D4").Copy()
$ws.Range("E4").PasteSpecial(-4122)

# F6: q4 (highlighted) -> "-" (plain), matching E6's style
$ws.Range("F6").Value = "-"
$ws.Range("E6").Copy()
$ws.Range("F6").PasteSpecial(-4122)

# G7: q5 (highlighted) -> "-" (plain), matching F7's style
$ws.Range("G7").Value = "-"
$ws.Range("F7").Copy()
$ws.Range("G7").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Update the view: scrolled right one column (top-left cell D1), selection moved to G7
$ws.Range("G7").Select()
$excel.ActiveWindow.ScrollColumn = 4
